$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (Fecha, Calidad, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
$ws.Range("D2").Value = 44447
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 21500
$ws.Range("S2").Value = 2150

$ws.Range("D3").Value = 44446
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21500
$ws.Range("S3").Value = 2150

$ws.Range("D4").Value = 44487
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 30
$ws.Range("N4").Value = 23000
$ws.Range("O4").Value = 24000
$ws.Range("P4").Value = 23500
$ws.Range("S4").Value = 2350

$ws.Range("D5").Value = 44460
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 31000
$ws.Range("O5").Value = 32000
$ws.Range("P5").Value = 31500
$ws.Range("S5").Value = 3150

$ws.Range("D6").Value = 44460
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 30000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 30000
$ws.Range("S6").Value = 3000

$ws.Range("D7").Value = 44448
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21500
$ws.Range("S7").Value = 2150
